$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description text for the sample/demo user row (I2).
$ws.Range("I2").Value = "Desc written here"

# Reflect the final on-screen view left after the edit: the window is
# scrolled one column to the right and the cursor ends up on C8.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C8").Select() | Out-Null
